$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are rotated among rows 26-29
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Capture the current ("before") values for the affected rows/columns
$rows = @(26, 27, 28, 29)
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row content mapping: row 26 gets old row 28's data, row 27 gets old row 29's data,
# row 28 gets old row 27's data, row 29 gets old row 26's data.
$mapping = @{
    26 = 28
    27 = 29
    28 = 27
    29 = 26
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $before[$src][$c]
    }
}
